$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the status of the first task (F2) from "Added" to "In progress"
$ws.Range("F2").Value = "In progress"

# Move the active selection to F2 (matching the author's edit)
$ws.Range("F2").Select()
